$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8. This shifts the existing data rows
# (old rows 8-65) down to become rows 9-66, preserving all of their
# values/styles, and leaves a fresh (mostly blank) row 8 that inherits
# the formatting of the row above it (so D8 keeps the date style).
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with this week's newest price record.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = 44831
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 100112040
$ws.Range("G8").Value = "Cilantro"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 7000
$ws.Range("N8").Value = "$/caja 36 atados"
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 194
$ws.Range("Q8").Value = 36
$ws.Range("R8").Value = "Hortaliza"

# Append a brand-new row 67 with the next day's price record.
$ws.Range("A67").Value = 5
$ws.Range("B67").Value = "Macroferia Regional de Talca"
$ws.Range("C67").Value = "Maule"
$ws.Range("D67").Value = 44832
$ws.Range("D67").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E67").Value = 7
$ws.Range("F67").Value = 100112040
$ws.Range("G67").Value = "Cilantro"
$ws.Range("H67").Value = "Sin especificar"
$ws.Range("I67").Value = "Primera"
$ws.Range("J67").Value = 150
$ws.Range("K67").Value = 7000
$ws.Range("L67").Value = 7000
$ws.Range("M67").Value = 7000
$ws.Range("N67").Value = "$/caja 36 atados"
$ws.Range("O67").Value = "Región del Maule"
$ws.Range("P67").Value = 194
$ws.Range("Q67").Value = 36
$ws.Range("R67").Value = "Hortaliza"
